$d = $word.ActiveDocument

# Locate the start of the "heres a cat:" paragraph and the start of the
# "heres the last line" paragraph that follows it later in the document.
# Between those two points sits: the "heres a cat:" paragraph, a blank
# paragraph, the paragraph holding the cat picture, and another blank
# paragraph -- all of which are being removed.
$catRange = $d.Content
$foundCat = $catRange.Find.Execute("heres a cat:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundCat) {
    throw "Could not find the 'heres a cat:' paragraph"
}
$catStart = $catRange.Paragraphs.Item(1).Range.Start

$lastLineRange = $d.Content
$foundLastLine = $lastLineRange.Find.Execute("heres the last line", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundLastLine) {
    throw "Could not find the 'heres the last line' paragraph"
}
$lastLineStart = $lastLineRange.Paragraphs.Item(1).Range.Start

# Delete everything from the start of "heres a cat:" up to (but not
# including) the start of "heres the last line" -- this removes the cat
# caption paragraph, the blank paragraph after it, the paragraph containing
# the cat picture, and the blank paragraph that followed the picture.
$d.Range($catStart, $lastLineStart).Delete()
